$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B, C, D, E columns remain text so numeric-looking strings are not
# reinterpreted as numbers/dates by Excel (matches original inlineStr text cells).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.390.12"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "2.072.06"
$ws.Range("E3").Value = "  +3.79%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Value = "328.26"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").Value = "0.5190"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("D8").Value = "0.4317"
$ws.Range("E8").Value = "  +4.40%  "

$ws.Range("D9").Value = "0.08665"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").Value = "46.17"
$ws.Range("E10").Value = "  +7.26%  "

$ws.Range("D11").Value = "1.150"
$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("D12").Value = "24.13"
$ws.Range("E12").Value = "  -2.18%  "

$ws.Range("D13").Value = "2.084.90"
$ws.Range("E13").Value = "  +4.70%  "

$ws.Range("D14").Value = "6.606"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "7.673"
$ws.Range("E15").Value = "  +3.25%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "95.16"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "0.00001112"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").Value = "  +1.62%  "

$ws.Range("D20").Value = "18.67"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "6.224"
$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("D23").Value = "30.416.69"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").Value = "12.23"
$ws.Range("E24").Value = "  +3.47%  "

$ws.Range("D25").Value = "2.301"
$ws.Range("E25").Value = "  +3.30%  "

$ws.Range("D26").Value = "2.329.91"
$ws.Range("E26").Value = "  +4.53%  "

$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "161.81"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "2.494"
$ws.Range("E29").Value = "  +3.97%  "

$ws.Range("D30").Value = "130.71"
$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("D31").Value = "1.169"
$ws.Range("E31").Value = "  +2.91%  "

$ws.Range("D32").Value = "0.1067"
$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("D33").Value = "6.014"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.840"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.498"
$ws.Range("E35").Value = "  +12.47%  "

$ws.Range("D36").Value = "0.02553"
$ws.Range("E36").Value = "  +1.31%  "

$ws.Range("D37").Value = "9.585"
$ws.Range("E37").Value = "  +6.38%  "

$ws.Range("D38").Value = "5.444"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").Value = "0.06593"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "12.40"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2235"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("D42").Value = "0.6697"
$ws.Range("E42").Value = "  +1.24%  "

$ws.Range("D43").Value = "1.239"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "1.005"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.84"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6289"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").Value = "2.185"
$ws.Range("E47").Value = "  -0.64%  "

$ws.Range("D48").Value = "3.621"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").Value = "1.225"
$ws.Range("E49").Value = "  -3.30%  "

$ws.Range("D50").Value = "81.57"
$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("D51").Value = "1.180"
$ws.Range("E51").Value = "  +6.81%  "

Write-Output "Update complete"